$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

# Column A must be stored as text (a date-like string), not auto-converted
# into a date serial number by Excel's type inference.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/15/2026"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 13219.41
$ws.Cells.Item($row, 3).Value = 0.222814822013932
$ws.Cells.Item($row, 4).Value = 0.777185177986068
$ws.Cells.Item($row, 5).Value = -114.52
$ws.Cells.Item($row, 6).Value = -18.03
$ws.Cells.Item($row, 7).Value = -20417.83
$ws.Cells.Item($row, 8).Value = -66.53
$ws.Cells.Item($row, 9).Value = -207.37
$ws.Cells.Item($row, 10).Value = -6.58
